$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    foreach ($r in 6..11) {
        $ws.Range("E$r").Copy($ws.Range("G$r"))
        $ws.Range("F$r").Copy($ws.Range("H$r"))
    }
    foreach ($r in 19..24) {
        $ws.Range("E$r").Copy($ws.Range("G$r"))
        $ws.Range("F$r").Copy($ws.Range("H$r"))
    }

    # Force the dependent IF() formulas to re-evaluate against the new
    # "Actual Output" values we just copied in (Copy doesn't dirty them).
    foreach ($r in (6..11) + (19..24)) {
        $cell = $ws.Range("I$r")
        $cell.Formula = $cell.Formula
    }

    $ws.Range("G19").Select()
}
